# Update the "Förändrad" (Changed) date column (C) for every data row
# from serial date 45186 (2023-09-17) to serial date 45188 (2023-09-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 319
}

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45188
